# feat: (upload-service) upload/download settings with pagination parameters
#
# Adds four new "pagination" settings columns (C:F) to the "settings" sheet,
# with their header labels in row 1 and values in row 2, and makes the
# "settings" sheet the active/selected tab (it was "language_English" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# --- New pagination columns: headers (row 1) ---------------------------------
$ws.Range("C1").Value = "languages pagination"
$ws.Range("D1").Value = "tags pagination"
$ws.Range("E1").Value = "notes pagination"
$ws.Range("F1").Value = "vocabulary pagination"

# --- New pagination columns: values (row 2) ----------------------------------
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2

# --- Column widths for the newly introduced columns --------------------------
$ws.Columns.Item(3).ColumnWidth = 17.506666666666668
$ws.Columns.Item(4).ColumnWidth = 14.856666666666666
$ws.Columns.Item(5).ColumnWidth = 15.696666666666667
$ws.Columns.Item(6).ColumnWidth = 15.976666666666665

# --- Make "settings" the active sheet / selected tab --------------------------
$ws.Activate()
$ws.Range("F2").Select()
